$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the replacement shipment's data (helpers for issue and
# declared sum). The digit-only fields are stored as text in the source
# data (numberStoredAsText), so they are entered with a leading apostrophe
# to keep Excel from auto-converting them to numbers.
$ws.Range("A2").Value = "'20230810"
$ws.Range("B2").Value = "'951690564664"
$ws.Range("C2").Value = 1559
$ws.Range("D2").Value = 250
$ws.Range("E2").Value = 250
$ws.Range("F2").Value = "'20230810"
$ws.Range("J2").Value = "'19749"
$ws.Range("L2").Value = "Абубякяров Ильдус"
$ws.Range("M2").Value = "'79299431323"
$ws.Range("AX2").Value = 1

# Remove the now-superseded rows 3 and 4
$ws.Range("A3:A4").EntireRow.Delete()
